$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_A_Relations")

# New rows 8 and 9 describing the predator/prey relationship between
# 密毛魔芋 (plant) and 綠背斜紋天蛾 (moth), matching List A & B update.
$ws.Range("A8").Value = "密毛魔芋"
$ws.Range("B8").Value = "綠背斜紋天蛾"
$ws.Range("C8").Value = "被吃"

$ws.Range("A9").Value = "綠背斜紋天蛾"
$ws.Range("B9").Value = "密毛魔芋"
$ws.Range("C9").Value = "吃"

# Match the formatting used by the other data rows (A2:C7 use the same cell
# style - left aligned, default border/font). Copy formats only so the
# values entered above are preserved.
$ws.Range("A2:C2").Copy()
$ws.Range("A8:C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false
